$d = $word.ActiveDocument

$d.Content.Find.Execute("821×5=", $true, $true, $false, $false, $false, $true, 1, $false, "432×9=", 2) | Out-Null
$d.Content.Find.Execute("393×8=", $true, $true, $false, $false, $false, $true, 1, $false, "108×8=", 2) | Out-Null
$d.Content.Find.Execute("228×8=", $true, $true, $false, $false, $false, $true, 1, $false, "937×3=", 2) | Out-Null
$d.Content.Find.Execute("916×5=", $true, $true, $false, $false, $false, $true, 1, $false, "597×3=", 2) | Out-Null
$d.Content.Find.Execute("450×7=", $true, $true, $false, $false, $false, $true, 1, $false, "863×8=", 2) | Out-Null
$d.Content.Find.Execute("934×4=", $true, $true, $false, $false, $false, $true, 1, $false, "829×4=", 2) | Out-Null
$d.Content.Find.Execute("431×3=", $true, $true, $false, $false, $false, $true, 1, $false, "195×2=", 2) | Out-Null
$d.Content.Find.Execute("385×3=", $true, $true, $false, $false, $false, $true, 1, $false, "258×5=", 2) | Out-Null
$d.Content.Find.Execute("212×9=", $true, $true, $false, $false, $false, $true, 1, $false, "504×7=", 2) | Out-Null
$d.Content.Find.Execute("150×6=", $true, $true, $false, $false, $false, $true, 1, $false, "103×9=", 2) | Out-Null
$d.Content.Find.Execute("994×5=", $true, $true, $false, $false, $false, $true, 1, $false, "908×4=", 2) | Out-Null
$d.Content.Find.Execute("662×2=", $true, $true, $false, $false, $false, $true, 1, $false, "877×2=", 2) | Out-Null
$d.Content.Find.Execute("839×3=", $true, $true, $false, $false, $false, $true, 1, $false, "170×2=", 2) | Out-Null
$d.Content.Find.Execute("386×3=", $true, $true, $false, $false, $false, $true, 1, $false, "361×9=", 2) | Out-Null
$d.Content.Find.Execute("510×9=", $true, $true, $false, $false, $false, $true, 1, $false, "815×4=", 2) | Out-Null
$d.Content.Find.Execute("761×2=", $true, $true, $false, $false, $false, $true, 1, $false, "197×9=", 2) | Out-Null
$d.Content.Find.Execute("837×8=", $true, $true, $false, $false, $false, $true, 1, $false, "938×5=", 2) | Out-Null
$d.Content.Find.Execute("115×5=", $true, $true, $false, $false, $false, $true, 1, $false, "703×5=", 2) | Out-Null
$d.Content.Find.Execute("773×5=", $true, $true, $false, $false, $false, $true, 1, $false, "514×3=", 2) | Out-Null
$d.Content.Find.Execute("697×6=", $true, $true, $false, $false, $false, $true, 1, $false, "157×6=", 2) | Out-Null
$d.Content.Find.Execute("225×9=", $true, $true, $false, $false, $false, $true, 1, $false, "509×2=", 2) | Out-Null
$d.Content.Find.Execute("501×3=", $true, $true, $false, $false, $false, $true, 1, $false, "622×2=", 2) | Out-Null
$d.Content.Find.Execute("589×5=", $true, $true, $false, $false, $false, $true, 1, $false, "829×2=", 2) | Out-Null
$d.Content.Find.Execute("935×2=", $true, $true, $false, $false, $false, $true, 1, $false, "383×7=", 2) | Out-Null
$d.Content.Find.Execute("102×2=", $true, $true, $false, $false, $false, $true, 1, $false, "837×7=", 2) | Out-Null
